$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial value that was bumped from
# 45172 (2023-09-03) to 45175 (2023-09-06) for every data row (rows 2-395).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 395 }

$ws.Range("C2:C$lastRow").Value = 45175
